# Solve the problem of sharing memory
# Adds a new "D&D first level" row (row 20) and a new results column (J)
# recording the "s x/y" success-ratio notes for several existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data in column J for existing rows (entered in the same order
#     the original author typed them, so that new shared strings line up) ---

# Row 20 is a brand-new level entry.
$ws.Range("A20").Value = "D&D first level"
$ws.Range("J20").Value = "4  mis"

$ws.Range("J14").Value = "s 5/6"

$ws.Range("J12").Value = "s  5/5 1mis`n3times tried"
$ws.Range("J12").WrapText = $true

$ws.Range("J10").Value = "s 2/4"

$ws.Range("J16").Value = "s 8/11"

$ws.Range("J15").Value = "s  8/11"

$ws.Range("J11").Value = "s"
$ws.Range("J17").Value = "s"
$ws.Range("J18").Value = "s"
$ws.Range("I20").Value = "s"

# --- Update the view to match the author's saved selection ---
$ws.Range("O6").Select()
